# Experiment order generation script re-run: regenerates randomized
# per-task stimulus-file orders (new timestamps) for each task-order sheet.
# Sheet 1 -> now vSAT task order (SAT/vSAT stims)
# Sheet 2 -> now TOL task order (MM/ZM stims)
# Sheet 3 -> now NB task order (OB/ZB/TB stims)
# Sheet 4 -> now GNG task order (go/GNG stims)
# Sheet 5 -> now RS task order (eyes closed/open) -- unchanged content

$wb = $excel.ActiveWorkbook

# --- Sheet 1: rename + new vSAT task order content ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "vSAT_TO-1651588976680924"
$ws1.Range("B2").Value = "SAT_stims-16515889766340487.csv"
$ws1.Range("B3").Value = "vSAT_stims-16515889766652992.csv"
$ws1.Range("B4").Value = "SAT_stims-16515889766184244.csv"
$ws1.Range("B5").Value = "vSAT_stims-16515889766496737.csv"

# --- Sheet 2: rename + new TOL task order content ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "TOL_TO-1651588976727799"
$ws2.Range("B2").Value = "MM_stims-16515889766966326.csv"
$ws2.Range("B3").Value = "ZM_stims-1651588976680924.csv"
$ws2.Range("B4").Value = "MM_stims-1651588976712181.csv"
$ws2.Range("B5").Value = "ZM_stims-16515889766966326.csv"
$ws2.Range("B6").Value = "MM_stims-1651588976727799.csv"
$ws2.Range("B7").Value = "ZM_stims-1651588976712181.csv"

# --- Sheet 3: rename + new NB task order content ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "NB_TO-1651588978221365"
$ws3.Range("B2").Value = "TB-16515889774343798.csv"
$ws3.Range("B3").Value = "ZB-match_3-1651588977078328.csv"
$ws3.Range("B4").Value = "TB-16515889782057402.csv"
$ws3.Range("B5").Value = "TB-1651588977872623.csv"
$ws3.Range("B6").Value = "OB-16515889771899147.csv"
$ws3.Range("B7").Value = "ZB-match_4-1651588976871081.csv"
$ws3.Range("B8").Value = "OB-16515889773049033.csv"
$ws3.Range("B9").Value = "OB-16515889770939543.csv"
$ws3.Range("B10").Value = "ZB-match_0-16515889769510245.csv"

# --- Sheet 4: rename + new GNG task order content ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "GNG_TO-16515889782608345"
$ws4.Range("B2").Value = "go_stims-1651588978221365.csv"
$ws4.Range("B3").Value = "GNG_stims-16515889782360344.csv"
$ws4.Range("B4").Value = "go_stims-16515889782370374.csv"
$ws4.Range("B5").Value = "GNG_stims-16515889782608345.csv"

# --- Sheet 5: rename only (RS task order content -- eyes closed/open -- unchanged) ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "RS_TO-16515889782608345"
